$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-13 Tuesday" "2024-02-14 Wednesday"

Replace-Text "625×5=" "611×8="
Replace-Text "815×8=" "539×2="
Replace-Text "223×3=" "212×8="
Replace-Text "590×3=" "711×6="
Replace-Text "508×2=" "818×8="
Replace-Text "393×9=" "470×8="
Replace-Text "756×5=" "947×5="
Replace-Text "566×4=" "715×3="
Replace-Text "404×7=" "725×6="
Replace-Text "456×5=" "692×4="
Replace-Text "983×9=" "305×2="
Replace-Text "553×8=" "443×9="
Replace-Text "828×5=" "420×5="
Replace-Text "791×5=" "885×9="
Replace-Text "942×6=" "984×2="
Replace-Text "273×2=" "529×7="
Replace-Text "500×2=" "340×2="
Replace-Text "647×9=" "309×4="
Replace-Text "263×4=" "107×2="
Replace-Text "803×8=" "436×7="
Replace-Text "865×7=" "768×7="
Replace-Text "679×7=" "628×4="
Replace-Text "874×5=" "852×6="
Replace-Text "638×4=" "132×9="
Replace-Text "449×9=" "185×5="
